$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '43.861.83'
Set-TextValue 'E2' '  -0.46%  '
Set-TextValue 'D3' '2.347.30'
Set-TextValue 'E3' '  -0.50%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '239.70'
Set-TextValue 'E5' '  -0.87%  '
Set-TextValue 'D6' '0.666'
Set-TextValue 'E6' '  -4.08%  '
Set-TextValue 'D7' '72.90'
Set-TextValue 'E7' '  -4.84%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'E9' '  -5.44%  '
Set-TextValue 'E10' '  -1.17%  '
Set-TextValue 'D11' '59.51'
Set-TextValue 'E11' '  +3.74%  '
Set-TextValue 'D12' '32.86'
Set-TextValue 'E12' '  -1.44%  '
Set-TextValue 'D14' '7.25'
Set-TextValue 'E14' '  -3.17%  '
Set-TextValue 'D15' '2.696.59'
Set-TextValue 'E15' '  -0.31%  '
Set-TextValue 'D16' '16.10'
Set-TextValue 'E16' '  -3.90%  '
Set-TextValue 'E17' '  -2.71%  '
Set-TextValue 'D18' '2.343.85'
Set-TextValue 'E18' '  -0.48%  '
Set-TextValue 'D19' '43.854.38'
Set-TextValue 'E19' '  -0.26%  '
Set-TextValue 'D20' '0.0000102'
Set-TextValue 'E20' '  -0.79%  '
Set-TextValue 'D21' '6.68'
Set-TextValue 'E21' '  -0.32%  '
Set-TextValue 'D22' '78.59'
Set-TextValue 'E22' '  +0.85%  '
Set-TextValue 'D23' '251.77'
Set-TextValue 'E23' '  -3.73%  '
Set-TextValue 'E24' '  +0.21%  '
Set-TextValue 'E25' '  +2.97%  '
Set-TextValue 'E26' '  +1.03%  '
Set-TextValue 'E27' '  -1.78%  '
Set-TextValue 'B28' 'Cosmos'
Set-TextValue 'C28' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D28' '10.41'
Set-TextValue 'E28' '  -4.95%  '
Set-TextValue 'B29' 'Toncoin'
Set-TextValue 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '2.32'
Set-TextValue 'E29' '  +0.82%  '
Set-TextValue 'D30' '176.53'
Set-TextValue 'E30' '  +0.81%  '
Set-TextValue 'E31' '  -4.33%  '
Set-TextValue 'E32' '  -1.35%  '
Set-TextValue 'E33' '  -2.74%  '
Set-TextValue 'E34' '  -2.43%  '
Set-TextValue 'E35' '  -6.17%  '
Set-TextValue 'E36' '  -2.41%  '
Set-TextValue 'E37' '  -1.68%  '
Set-TextValue 'E38' '  -0.34%  '
Set-TextValue 'E39' '  -2.54%  '
Set-TextValue 'D40' '5.68'
Set-TextValue 'E40' '  +20.10%  '
Set-TextValue 'E41' '  -4.58%  '
Set-TextValue 'D42' '65.08'
Set-TextValue 'E42' '  +15.06%  '
Set-TextValue 'E43' '  +0.62%  '
Set-TextValue 'E44' '  -2.47%  '
Set-TextValue 'D45' '18.80'
Set-TextValue 'E45' '  -2.63%  '
Set-TextValue 'E46' '  -10.91%  '
Set-TextValue 'E47' '  +0.10%  '
Set-TextValue 'E48' '  -3.19%  '
Set-TextValue 'E49' '  -3.28%  '
Set-TextValue 'D50' '98.15'
Set-TextValue 'E50' '  -4.03%  '
Set-TextValue 'E51' '  -5.15%  '
